$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New vendor-selection rows appended below the existing client_id/device_name/
# selected_vendor/selected_model/created_at/updated_at records (rows 4-9).
# Columns: A=client_id, B=device_name, C=selected_vendor, D=selected_model,
#          E=created_at, F=updated_at
$newRows = @(
    @(86, "PV_RPT_DB.ahm.lambdacro.com", "Microsoft", "Windows Server 2022", "2025-12-25 17:25:46", "2025-12-25 17:25:46"),
    @(86, "AHM_PRD_NAS",                 "Synology",  "NAS",                 "2025-12-25 17:25:49", "2025-12-25 17:25:49"),
    @(86, "MEH_PRD_NAS",                 "Synology",  "NAS",                 "2025-12-25 17:25:49", "2025-12-25 17:25:49"),
    @(86, "HOST-3.123.68.65",            "Fortinet",  "FortiWeb Cloud",      "2025-12-25 17:25:49", "2025-12-25 17:25:49"),
    @(86, "AHM_VCenter",                 "VMware",    "vCenter",             "2025-12-25 17:26:41", "2025-12-25 17:26:41"),
    @(86, "LTRPDC.ahm.lambdacro.com",    "Microsoft", "Windows",             "2025-12-25 17:26:53", "2025-12-25 17:26:53")
)

$startRow = 4
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
    $ws.Cells.Item($row, 4).Value = $rowValues[3]
    $ws.Cells.Item($row, 5).Value = $rowValues[4]
    $ws.Cells.Item($row, 6).Value = $rowValues[5]
}
